$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(273).Insert()

$ws.Cells.Item(273, 1).Value = 7
$ws.Cells.Item(273, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(273, 3).Value = "Ñuble"
$ws.Cells.Item(273, 4).Value = 44855
$ws.Cells.Item(273, 5).Value = 16
$ws.Cells.Item(273, 6).Value = 100114013
$ws.Cells.Item(273, 7).Value = "Zanahoria"
$ws.Cells.Item(273, 8).Value = "Sin especificar"
$ws.Cells.Item(273, 9).Value = "Primera"
$ws.Cells.Item(273, 10).Value = 80
$ws.Cells.Item(273, 11).Value = 20000
$ws.Cells.Item(273, 12).Value = 20000
$ws.Cells.Item(273, 13).Value = 20000
$ws.Cells.Item(273, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(273, 15).Value = "Región de Ñuble"
$ws.Cells.Item(273, 16).Value = 1000
$ws.Cells.Item(273, 17).Value = 20
$ws.Cells.Item(273, 18).Value = "Hortaliza"
